# Updated symbol list on Sat Jan 21 20:06:26 UTC 2023 with GitHub Actions
# Refreshes price (D), 1h volume change (E) and collection hour (G) for each
# coin row on Sheet1, keeping every touched cell text-typed (as the source
# data always was) instead of letting Excel auto-coerce the numeric-looking
# strings into Number/Percentage cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "304.07" },
    @{ Cell = "E2"; Value = "3.07%" },
    @{ Cell = "G2"; Value = "20" },
    @{ Cell = "D3"; Value = "35.63" },
    @{ Cell = "E3"; Value = "13.69%" },
    @{ Cell = "G3"; Value = "20" },
    @{ Cell = "D4"; Value = "5.084" },
    @{ Cell = "E4"; Value = "2.07%" },
    @{ Cell = "G4"; Value = "20" },
    @{ Cell = "D5"; Value = "0.07826" },
    @{ Cell = "E5"; Value = "2.82%" },
    @{ Cell = "G5"; Value = "20" },
    @{ Cell = "D6"; Value = "2.298" },
    @{ Cell = "E6"; Value = "2.40%" },
    @{ Cell = "G6"; Value = "20" },
    @{ Cell = "D7"; Value = "8.118" },
    @{ Cell = "E7"; Value = "4.11%" },
    @{ Cell = "G7"; Value = "20" },
    @{ Cell = "D8"; Value = "4.018" },
    @{ Cell = "E8"; Value = "6.33%" },
    @{ Cell = "G8"; Value = "20" },
    @{ Cell = "D9"; Value = "0.9278" },
    @{ Cell = "E9"; Value = "0.59%" },
    @{ Cell = "G9"; Value = "20" },
    @{ Cell = "D10"; Value = "0.09701" },
    @{ Cell = "E10"; Value = "1.85%" },
    @{ Cell = "G10"; Value = "20" },
    @{ Cell = "D11"; Value = "0.1827" },
    @{ Cell = "E11"; Value = "4.64%" },
    @{ Cell = "G11"; Value = "20" },
    @{ Cell = "D12"; Value = "0.08720" },
    @{ Cell = "E12"; Value = "3.86%" },
    @{ Cell = "G12"; Value = "20" },
    @{ Cell = "D13"; Value = "0.03420" },
    @{ Cell = "E13"; Value = "4.37%" },
    @{ Cell = "G13"; Value = "20" },
    @{ Cell = "D14"; Value = "0.09948" },
    @{ Cell = "E14"; Value = "-0.09%" },
    @{ Cell = "G14"; Value = "20" },
    @{ Cell = "D15"; Value = "0.001484" },
    @{ Cell = "E15"; Value = "-0.79%" },
    @{ Cell = "G15"; Value = "20" },
    @{ Cell = "D16"; Value = "0.005681" },
    @{ Cell = "E16"; Value = "-0.57%" },
    @{ Cell = "G16"; Value = "20" },
    @{ Cell = "D17"; Value = "3.489" },
    @{ Cell = "E17"; Value = "0.33%" },
    @{ Cell = "G17"; Value = "20" },
    @{ Cell = "G18"; Value = "20" },
    @{ Cell = "E19"; Value = "2.09%" },
    @{ Cell = "G19"; Value = "20" },
    @{ Cell = "E20"; Value = "0.29%" },
    @{ Cell = "G20"; Value = "20" },
    @{ Cell = "D21"; Value = "4.553" },
    @{ Cell = "E21"; Value = "11.48%" },
    @{ Cell = "G21"; Value = "20" },
    @{ Cell = "D22"; Value = "0.2235" },
    @{ Cell = "E22"; Value = "-2.41%" },
    @{ Cell = "G22"; Value = "20" },
    @{ Cell = "D23"; Value = "0.04680" },
    @{ Cell = "E23"; Value = "3.58%" },
    @{ Cell = "G23"; Value = "20" },
    @{ Cell = "D24"; Value = "0.001241" },
    @{ Cell = "E24"; Value = "1.92%" },
    @{ Cell = "G24"; Value = "20" },
    @{ Cell = "D25"; Value = "0.004547" },
    @{ Cell = "E25"; Value = "6.02%" },
    @{ Cell = "G25"; Value = "20" },
    @{ Cell = "D26"; Value = "0.0001300" },
    @{ Cell = "E26"; Value = "0.65%" },
    @{ Cell = "G26"; Value = "20" },
    @{ Cell = "D27"; Value = "0.0002699" },
    @{ Cell = "E27"; Value = "-20.36%" },
    @{ Cell = "G27"; Value = "20" },
    @{ Cell = "G28"; Value = "20" },
    @{ Cell = "G29"; Value = "20" },
    @{ Cell = "G30"; Value = "20" },
    @{ Cell = "G31"; Value = "20" },
    @{ Cell = "G32"; Value = "20" },
    @{ Cell = "G33"; Value = "20" },
    @{ Cell = "G34"; Value = "20" },
    @{ Cell = "G35"; Value = "20" },
    @{ Cell = "G36"; Value = "20" },
    @{ Cell = "G37"; Value = "20" },
    @{ Cell = "G38"; Value = "20" },
    @{ Cell = "E39"; Value = "6.13%" },
    @{ Cell = "G39"; Value = "20" },
    @{ Cell = "D40"; Value = "0.04712" },
    @{ Cell = "E40"; Value = "2.31%" },
    @{ Cell = "G40"; Value = "20" },
    @{ Cell = "D41"; Value = "0.008011" },
    @{ Cell = "E41"; Value = "7.40%" },
    @{ Cell = "G41"; Value = "20" },
    @{ Cell = "E42"; Value = "3.99%" },
    @{ Cell = "G42"; Value = "20" },
    @{ Cell = "D43"; Value = "0.008014" },
    @{ Cell = "E43"; Value = "-18.44%" },
    @{ Cell = "G43"; Value = "20" },
    @{ Cell = "D44"; Value = "0.002300" },
    @{ Cell = "E44"; Value = "7.14%" },
    @{ Cell = "G44"; Value = "20" },
    @{ Cell = "D45"; Value = "0.009126" },
    @{ Cell = "E45"; Value = "-2.03%" },
    @{ Cell = "G45"; Value = "20" },
    @{ Cell = "D46"; Value = "0.00006225" },
    @{ Cell = "E46"; Value = "2.15%" },
    @{ Cell = "G46"; Value = "20" },
    @{ Cell = "D47"; Value = "0.00000000750" },
    @{ Cell = "E47"; Value = "0.16%" },
    @{ Cell = "G47"; Value = "20" },
    @{ Cell = "D48"; Value = "5.721" },
    @{ Cell = "E48"; Value = "115.53%" },
    @{ Cell = "G48"; Value = "20" },
    @{ Cell = "D49"; Value = "0.002690" },
    @{ Cell = "E49"; Value = "34.72%" },
    @{ Cell = "G49"; Value = "20" },
    @{ Cell = "D50"; Value = "0.00002100" },
    @{ Cell = "E50"; Value = "0.16%" },
    @{ Cell = "G50"; Value = "20" },
    @{ Cell = "D51"; Value = "0.0002000" },
    @{ Cell = "E51"; Value = "0.16%" },
    @{ Cell = "G51"; Value = "20" }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    # Force text storage so "304.07" / "3.07%" / "20" stay literal strings
    # (matching the original inlineStr cells) rather than being parsed into
    # Number/Percentage values with an auto-assigned number format.
    $range.NumberFormat = "@"
    $range.Value = $u.Value
    # Drop the temporary "@" text format again so the cell's style index
    # reverts to the original default (no explicit style), exactly like the
    # untouched cells around it.
    $range.ClearFormats()
}
